$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Add the new "metadata" worksheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# --- Header row (B1:G1) ---
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Match the bold/bordered/centered header style used on the "data" sheet
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# --- Data row (row 2) ---
$metaSheet.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("B2").Value = "Familial Hirschsprung Disease"
$metaSheet.Range("C2").Value = 63

$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.9"

$metaSheet.Range("E2").Value = "2021-01-19T18:07:25.634925Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:14.737262"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/63/?format=json"

$excel.CutCopyMode = $false

[void]$dataSheet.Activate()
[void]$dataSheet.Range("A1").Select()

# --- Update column F timestamps on the "data" sheet ---
$dataSheet.Cells.Item(2, 6).Value = "2021-10-05 14:20:14.740900"
$dataSheet.Cells.Item(3, 6).Value = "2021-10-05 14:20:14.740907"
$dataSheet.Cells.Item(4, 6).Value = "2021-10-05 14:20:14.740911"
$dataSheet.Cells.Item(5, 6).Value = "2021-10-05 14:20:14.740913"
$dataSheet.Cells.Item(6, 6).Value = "2021-10-05 14:20:14.740916"
$dataSheet.Cells.Item(7, 6).Value = "2021-10-05 14:20:14.740919"
$dataSheet.Cells.Item(8, 6).Value = "2021-10-05 14:20:14.740921"
$dataSheet.Cells.Item(9, 6).Value = "2021-10-05 14:20:14.740924"
$dataSheet.Cells.Item(10, 6).Value = "2021-10-05 14:20:14.740927"
$dataSheet.Cells.Item(11, 6).Value = "2021-10-05 14:20:14.740929"
$dataSheet.Cells.Item(12, 6).Value = "2021-10-05 14:20:14.740932"
$dataSheet.Cells.Item(13, 6).Value = "2021-10-05 14:20:14.740934"
$dataSheet.Cells.Item(14, 6).Value = "2021-10-05 14:20:14.740937"
$dataSheet.Cells.Item(15, 6).Value = "2021-10-05 14:20:14.740939"
$dataSheet.Cells.Item(16, 6).Value = "2021-10-05 14:20:14.740942"
$dataSheet.Cells.Item(17, 6).Value = "2021-10-05 14:20:14.740944"
$dataSheet.Cells.Item(18, 6).Value = "2021-10-05 14:20:14.740947"
$dataSheet.Cells.Item(19, 6).Value = "2021-10-05 14:20:14.740950"
$dataSheet.Cells.Item(20, 6).Value = "2021-10-05 14:20:14.740952"
$dataSheet.Cells.Item(21, 6).Value = "2021-10-05 14:20:14.740955"
$dataSheet.Cells.Item(22, 6).Value = "2021-10-05 14:20:14.740958"
$dataSheet.Cells.Item(23, 6).Value = "2021-10-05 14:20:14.740960"
$dataSheet.Cells.Item(24, 6).Value = "2021-10-05 14:20:14.740963"
$dataSheet.Cells.Item(25, 6).Value = "2021-10-05 14:20:14.740966"
$dataSheet.Cells.Item(26, 6).Value = "2021-10-05 14:20:14.740968"
$dataSheet.Cells.Item(27, 6).Value = "2021-10-05 14:20:14.740971"
$dataSheet.Cells.Item(28, 6).Value = "2021-10-05 14:20:14.740974"
$dataSheet.Cells.Item(29, 6).Value = "2021-10-05 14:20:14.740976"
$dataSheet.Cells.Item(30, 6).Value = "2021-10-05 14:20:14.740979"
$dataSheet.Cells.Item(31, 6).Value = "2021-10-05 14:20:14.740981"
$dataSheet.Cells.Item(32, 6).Value = "2021-10-05 14:20:14.740984"
$dataSheet.Cells.Item(33, 6).Value = "2021-10-05 14:20:14.740986"
$dataSheet.Cells.Item(34, 6).Value = "2021-10-05 14:20:14.740989"
$dataSheet.Cells.Item(35, 6).Value = "2021-10-05 14:20:14.740992"
$dataSheet.Cells.Item(36, 6).Value = "2021-10-05 14:20:14.740994"
$dataSheet.Cells.Item(37, 6).Value = "2021-10-05 14:20:14.740997"
$dataSheet.Cells.Item(38, 6).Value = "2021-10-05 14:20:14.741000"
$dataSheet.Cells.Item(39, 6).Value = "2021-10-05 14:20:14.741002"
$dataSheet.Cells.Item(40, 6).Value = "2021-10-05 14:20:14.741005"
$dataSheet.Cells.Item(41, 6).Value = "2021-10-05 14:20:14.741007"
$dataSheet.Cells.Item(42, 6).Value = "2021-10-05 14:20:14.741010"
$dataSheet.Cells.Item(43, 6).Value = "2021-10-05 14:20:14.741013"
$dataSheet.Cells.Item(44, 6).Value = "2021-10-05 14:20:14.741015"
$dataSheet.Cells.Item(45, 6).Value = "2021-10-05 14:20:14.741018"
$dataSheet.Cells.Item(46, 6).Value = "2021-10-05 14:20:14.741020"
$dataSheet.Cells.Item(47, 6).Value = "2021-10-05 14:20:14.741023"
$dataSheet.Cells.Item(48, 6).Value = "2021-10-05 14:20:14.741026"
$dataSheet.Cells.Item(49, 6).Value = "2021-10-05 14:20:14.741028"
$dataSheet.Cells.Item(50, 6).Value = "2021-10-05 14:20:14.741031"
$dataSheet.Cells.Item(51, 6).Value = "2021-10-05 14:20:14.741033"
$dataSheet.Cells.Item(52, 6).Value = "2021-10-05 14:20:14.741036"
$dataSheet.Cells.Item(53, 6).Value = "2021-10-05 14:20:14.741038"
$dataSheet.Cells.Item(54, 6).Value = "2021-10-05 14:20:14.741041"
$dataSheet.Cells.Item(55, 6).Value = "2021-10-05 14:20:14.741044"
$dataSheet.Cells.Item(56, 6).Value = "2021-10-05 14:20:14.741046"
$dataSheet.Cells.Item(57, 6).Value = "2021-10-05 14:20:14.741049"
$dataSheet.Cells.Item(58, 6).Value = "2021-10-05 14:20:14.741051"
$dataSheet.Cells.Item(59, 6).Value = "2021-10-05 14:20:14.741054"
$dataSheet.Cells.Item(60, 6).Value = "2021-10-05 14:20:14.741056"
$dataSheet.Cells.Item(61, 6).Value = "2021-10-05 14:20:14.741059"
$dataSheet.Cells.Item(62, 6).Value = "2021-10-05 14:20:14.741061"
$dataSheet.Cells.Item(63, 6).Value = "2021-10-05 14:20:14.741064"
$dataSheet.Cells.Item(64, 6).Value = "2021-10-05 14:20:14.741066"
